# Add target and sales sigh in columns name
# Rewrites the numeric values across columns A:L for rows 2-52 (swapping which
# columns hold the "target"/"sales" figures for several rows), then removes the
# now-unused trailing rows 53-61 so the sheet dimension shrinks to A1:L52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("K2").Value = 0.1735427119226257
$ws.Range("L2").Value = 0.4191043763702654
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("K3").Value = 0.2644460372154296
$ws.Range("L3").Value = 1.35402952365778
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0.385638222623379
$ws.Range("L4").Value = 0.9778791529581891
$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("K5").Value = 2.97488563686306
$ws.Range("L5").Value = 1.289494588390176
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("K6").Value = 5.164777939627022
$ws.Range("L6").Value = 6.98482447092498
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("K7").Value = 7.272266023424316
$ws.Range("L7").Value = 4.513431745525934
$ws.Range("C8").Value = 0
$ws.Range("K8").Value = 16.96767760100764
$ws.Range("L8").Value = 16.33341654641623
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 27.68296356525348
$ws.Range("L9").Value = 25.14514447360845
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("K10").Value = 39.11380226206305
$ws.Range("L10").Value = 42.98267512214798
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("I11").Value = 0.6226210791404801
$ws.Range("J11").Value = 1.865849732498522
$ws.Range("I12").Value = 1.639316513031779
$ws.Range("J12").Value = 5.038613462025728
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("I13").Value = 4.571323494546493
$ws.Range("J13").Value = 1.968274023369803
$ws.Range("I14").Value = 8.27556144976719
$ws.Range("J14").Value = 8.916870955805498
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("I15").Value = 8.511835740741928
$ws.Range("J15").Value = 7.377389821827624
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("I16").Value = 11.15996241563942
$ws.Range("J16").Value = 8.109497208726548
$ws.Range("A17").Value = 0
$ws.Range("B17").Value = 0
$ws.Range("I17").Value = 20.64977164777899
$ws.Range("J17").Value = 14.95888257761051
$ws.Range("A18").Value = 0
$ws.Range("B18").Value = 0
$ws.Range("I18").Value = 21.87112668867043
$ws.Range("J18").Value = 18.69835379090232
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("I19").Value = 22.69848097068329
$ws.Range("J19").Value = 33.06626842723345
$ws.Range("G20").Value = 0.4222366318254477
$ws.Range("H20").Value = 0.9923481590328036
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0.7037746599251754
$ws.Range("H21").Value = 2.366312778236024
$ws.Range("G22").Value = 4.644809487030747
$ws.Range("H22").Value = 3.549390253434932
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("G23").Value = 6.896991667266719
$ws.Range("H23").Value = 5.657065081836256
$ws.Range("C24").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("G24").Value = 87.33218755395191
$ws.Range("H24").Value = 87.43488372745999
$ws.Range("C25").Value = 0
$ws.Range("E25").Value = 0.03809864712110639
$ws.Range("F25").Value = 0.3567609425795917
$ws.Range("E26").Value = 0.08890249003626782
$ws.Range("F26").Value = 0.4252183902787009
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("E27").Value = 0.09524238273750128
$ws.Range("F27").Value = 0.4435232593774707
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("E28").Value = 0.1600007656997999
$ws.Range("F28").Value = 0.5242184089857995
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0.3714491042349933
$ws.Range("F29").Value = 0.2882930768633495
$ws.Range("E30").Value = 0.9870412508485481
$ws.Range("F30").Value = 0.8797336972882004
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("A31").Value = 0
$ws.Range("B31").Value = 0
$ws.Range("E31").Value = 98.25926535932177
$ws.Range("F31").Value = 97.08225222462688
$ws.Range("C32").Value = 0.2992388439095395
$ws.Range("D32").Value = 0.9625724567901734
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("C33").Value = 0.4103827450376586
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("A34").Value = 0
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = 0.7694419851931689
$ws.Range("D34").Value = 0.6818435200903089
$ws.Range("C35").Value = 1.346523474088045
$ws.Range("D35").Value = 1.286820936254987
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("A36").Value = 0
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 1.795444468680544
$ws.Range("D36").Value = 1.143891696724708
$ws.Range("C37").Value = 2.564806617310563
$ws.Range("D37").Value = 4.329706352573461
$ws.Range("A38").Value = 0
$ws.Range("B38").Value = 0
$ws.Range("C38").Value = 3.590888937361088
$ws.Range("D38").Value = 0.8864359925428781
$ws.Range("C39").Value = 8.771985349990661
$ws.Range("D39").Value = 6.176984643614592
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("C40").Value = 8.848976309925538
$ws.Range("D40").Value = 28.5972924181177
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("C41").Value = 10.62338384421605
$ws.Range("D41").Value = 8.75099222169267
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("C42").Value = 14.42618183062388
$ws.Range("D42").Value = 16.29876543376681
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("C43").Value = 16.92847641898799
$ws.Range("D43").Value = 10.86697111888473
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("C44").Value = 29.62426917467527
$ws.Range("D44").Value = 20.01772320894699
$ws.Range("A45").Value = 0.2413350402456683
$ws.Range("B45").Value = 0.3675400741575188
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("A46").Value = 0.3619971947444142
$ws.Range("B46").Value = 0.3445637123025586
$ws.Range("A47").Value = 5.140664932818899
$ws.Range("B47").Value = 3.275697813628217
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("A48").Value = 6.033429662382592
$ws.Range("B48").Value = 8.04001062327734
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 0
$ws.Range("A49").Value = 6.083383622645102
$ws.Range("B49").Value = 3.452021567988243
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("A50").Value = 15.96549495923762
$ws.Range("B50").Value = 8.990112509716967
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("A51").Value = 29.36190764247839
$ws.Range("B51").Value = 28.90141678061865
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("A52").Value = 36.81178694544733
$ws.Range("B52").Value = 46.62863691831051
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 0

# Remove rows 53:61 entirely (data moved up into rows 2-52 above); this also
# shrinks the sheet dimension from A1:L61 to A1:L52 automatically.
$ws.Rows("53:61").Delete()
